$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 54, shifting existing rows 54-66 down to 56-68.
$ws.Range("A54:T55").EntireRow.Insert()

# Fill in new row 54 (Packham's Triumph)
$ws.Cells.Item(54, 1).Value = 1
$ws.Cells.Item(54, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(54, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(54, 4).Value = 45093
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 15
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100104
$ws.Cells.Item(54, 8).Value = "Frutos de pepita"
$ws.Cells.Item(54, 9).Value = 100104005
$ws.Cells.Item(54, 10).Value = "Pera"
$ws.Cells.Item(54, 11).Value = "Packham's Triumph"
$ws.Cells.Item(54, 12).Value = "Segunda"
$ws.Cells.Item(54, 13).Value = 300
$ws.Cells.Item(54, 14).Value = 19000
$ws.Cells.Item(54, 15).Value = 20000
$ws.Cells.Item(54, 16).Value = 19500
$ws.Cells.Item(54, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(54, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(54, 19).Value = 1083
$ws.Cells.Item(54, 20).Value = 18

# Fill in new row 55 (Winter Nelis)
$ws.Cells.Item(55, 1).Value = 1
$ws.Cells.Item(55, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(55, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(55, 4).Value = 45093
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 5).Value = 15
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100104
$ws.Cells.Item(55, 8).Value = "Frutos de pepita"
$ws.Cells.Item(55, 9).Value = 100104005
$ws.Cells.Item(55, 10).Value = "Pera"
$ws.Cells.Item(55, 11).Value = "Winter Nelis"
$ws.Cells.Item(55, 12).Value = "Segunda"
$ws.Cells.Item(55, 13).Value = 300
$ws.Cells.Item(55, 14).Value = 19000
$ws.Cells.Item(55, 15).Value = 20000
$ws.Cells.Item(55, 16).Value = 19500
$ws.Cells.Item(55, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(55, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(55, 19).Value = 1083
$ws.Cells.Item(55, 20).Value = 18
